$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 83.91225566666667
$ws.Range("N2").Value = 251.736767
$ws.Range("O2").Value = 0.9556261553553385
$ws.Range("P2").Value = 0.9556261553553385
$ws.Range("Q2").Value = 13.839228765825
$ws.Range("R2").Value = 124.553058892425
$ws.Range("S2").Value = 0.9556261553553385
$ws.Range("T2").Value = 0.9556261553553385

# Row 3
$ws.Range("O3").Value = 0.00439999103960854
$ws.Range("P3").Value = 0.00439999103960854
$ws.Range("S3").Value = 0.00439999103960854
$ws.Range("T3").Value = 0.00439999103960854

# Row 4
$ws.Range("M4").Value = 3.510050666666667
$ws.Range("N4").Value = 10.530152
$ws.Range("O4").Value = 0.03997385360505296
$ws.Range("P4").Value = 0.03997385360505297
$ws.Range("Q4").Value = 0.5788951062000001
$ws.Range("R4").Value = 5.210055955800001
$ws.Range("S4").Value = 0.03997385360505296
$ws.Range("T4").Value = 0.03997385360505297

$wb.Save()
